$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet references (by name, resolved fresh each time to avoid stale handles
# after sheets get moved/reordered)
# ---------------------------------------------------------------------------
$wsUser        = $wb.Worksheets.Item("User")
$wsDaily       = $wb.Worksheets.Item("Daily level")
$wsGroupMember = $wb.Worksheets.Item("Group Member")
$wsTarget      = $wb.Worksheets.Item("Target")

# ---------------------------------------------------------------------------
# "User" sheet: highlight the "Created Date" header (D1) like the other
# highlighted headers, then leave the active selection at F14.
# ---------------------------------------------------------------------------
$wsUser.Range("A1").Copy() | Out-Null
$wsUser.Range("D1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# "Daily level" sheet: the measurement dates were re-entered with new
# values. Some of the re-typed cells picked up a slightly different (but
# visually equivalent) date number format.
# ---------------------------------------------------------------------------
$wsDaily.Range("A2").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A2").Value = (Get-Date -Year 2020 -Month 3 -Day 12 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A3").Value = (Get-Date -Year 2020 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A4").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A4").Value = (Get-Date -Year 2020 -Month 4 -Day 17 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A5").Value = (Get-Date -Year 2020 -Month 5 -Day 5 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A6").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A6").Value = (Get-Date -Year 2020 -Month 5 -Day 23 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A7").Value = (Get-Date -Year 2020 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A8").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A8").Value = (Get-Date -Year 2020 -Month 6 -Day 28 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A9").Value = (Get-Date -Year 2020 -Month 7 -Day 16 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A10").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A10").Value = (Get-Date -Year 2020 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 0)

$wsDaily.Range("A11").NumberFormat = "m/d/yyyy;@"
$wsDaily.Range("A11").Value = (Get-Date -Year 2020 -Month 8 -Day 21 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# "Group Member" sheet: highlight the "Date Invited" header (A1), and the
# pending group-member record in row 5 now has an actual "Date Left Group"
# instead of the "pending" placeholder text.
# ---------------------------------------------------------------------------
$wsUser.Range("A1").Copy() | Out-Null
$wsGroupMember.Range("A1").PasteSpecial(-4122) | Out-Null

$wsGroupMember.Range("A5").Copy() | Out-Null
$wsGroupMember.Range("D5").PasteSpecial(-4122) | Out-Null
$wsGroupMember.Range("D5").Value = (Get-Date -Year 2020 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Highlight the "Target Name" header (A1) on the "Target" sheet the same
# way, then reorder the tabs so "Target" comes before "Challenge".
# ---------------------------------------------------------------------------
$wsUser.Range("A1").Copy() | Out-Null
$wsTarget.Range("A1").PasteSpecial(-4122) | Out-Null

$wb.Worksheets.Item("Target").Move($wb.Worksheets.Item("Challenge")) | Out-Null

# ---------------------------------------------------------------------------
# Selections / active sheet & cell state
# ---------------------------------------------------------------------------
$wsUser.Range("F14").Select() | Out-Null
$wsDaily.Range("B1").Select() | Out-Null
$wsGroupMember.Range("D7").Select() | Out-Null

# "Challenge" ends up as the last (rightmost) tab and is the active sheet.
$wb.Worksheets.Item("Challenge").Activate() | Out-Null
$wb.Worksheets.Item("Challenge").Range("B1").Select() | Out-Null
